$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.564.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.251"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0908"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.591.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.560.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.536"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0476"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.428.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.538"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("E43").Value = "  +7.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.983"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.735.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "85.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.48%  "
